# Auto-generated edit script to update cryptos worksheet values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'294.64"
$ws.Range("E2").Value = "'2.11%"
$ws.Range("G2").Value = "'2"
$ws.Range("D3").Value = "'31.15"
$ws.Range("E3").Value = "'0.97%"
$ws.Range("G3").Value = "'2"
$ws.Range("D4").Value = "'4.933"
$ws.Range("E4").Value = "'-0.21%"
$ws.Range("G4").Value = "'2"
$ws.Range("D5").Value = "'0.07409"
$ws.Range("E5").Value = "'4.16%"
$ws.Range("G5").Value = "'2"
$ws.Range("D6").Value = "'2.141"
$ws.Range("E6").Value = "'18.19%"
$ws.Range("G6").Value = "'2"
$ws.Range("D7").Value = "'7.746"
$ws.Range("E7").Value = "'1.05%"
$ws.Range("G7").Value = "'2"
$ws.Range("D8").Value = "'3.747"
$ws.Range("E8").Value = "'-0.29%"
$ws.Range("G8").Value = "'2"
$ws.Range("D9").Value = "'0.9145"
$ws.Range("E9").Value = "'2.04%"
$ws.Range("G9").Value = "'2"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.08874"
$ws.Range("E10").Value = "'17.69%"
$ws.Range("G10").Value = "'2"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1698"
$ws.Range("E11").Value = "'3.28%"
$ws.Range("G11").Value = "'2"
$ws.Range("D12").Value = "'0.08385"
$ws.Range("E12").Value = "'3.59%"
$ws.Range("G12").Value = "'2"
$ws.Range("D13").Value = "'0.03144"
$ws.Range("E13").Value = "'2.47%"
$ws.Range("G13").Value = "'2"
$ws.Range("D14").Value = "'0.1008"
$ws.Range("E14").Value = "'0.81%"
$ws.Range("G14").Value = "'2"
$ws.Range("D15").Value = "'0.001516"
$ws.Range("E15").Value = "'1.72%"
$ws.Range("G15").Value = "'2"
$ws.Range("D16").Value = "'0.005847"
$ws.Range("E16").Value = "'0.94%"
$ws.Range("G16").Value = "'2"
$ws.Range("D17").Value = "'3.498"
$ws.Range("E17").Value = "'0.85%"
$ws.Range("G17").Value = "'2"
$ws.Range("D18").Value = "'2.075"
$ws.Range("E18").Value = "'-2.51%"
$ws.Range("G18").Value = "'2"
$ws.Range("E19").Value = "'1.73%"
$ws.Range("G19").Value = "'2"
$ws.Range("D20").Value = "'0.1293"
$ws.Range("E20").Value = "'-0.20%"
$ws.Range("G20").Value = "'2"
$ws.Range("D21").Value = "'3.983"
$ws.Range("E21").Value = "'-1.66%"
$ws.Range("G21").Value = "'2"
$ws.Range("E22").Value = "'4.81%"
$ws.Range("G22").Value = "'2"
$ws.Range("D23").Value = "'0.04554"
$ws.Range("G23").Value = "'2"
$ws.Range("D24").Value = "'0.001211"
$ws.Range("E24").Value = "'-0.13%"
$ws.Range("G24").Value = "'2"
$ws.Range("D25").Value = "'0.004651"
$ws.Range("E25").Value = "'16.79%"
$ws.Range("G25").Value = "'2"
$ws.Range("D26").Value = "'0.0001299"
$ws.Range("E26").Value = "'3.71%"
$ws.Range("G26").Value = "'2"
$ws.Range("D27").Value = "'0.0003392"
$ws.Range("G27").Value = "'2"
$ws.Range("G28").Value = "'2"
$ws.Range("G29").Value = "'2"
$ws.Range("G30").Value = "'2"
$ws.Range("G31").Value = "'2"
$ws.Range("G32").Value = "'2"
$ws.Range("G33").Value = "'2"
$ws.Range("G34").Value = "'2"
$ws.Range("G35").Value = "'2"
$ws.Range("G36").Value = "'2"
$ws.Range("G37").Value = "'2"
$ws.Range("G38").Value = "'2"
$ws.Range("D39").Value = "'0.01624"
$ws.Range("E39").Value = "'1.52%"
$ws.Range("G39").Value = "'2"
$ws.Range("D40").Value = "'0.04467"
$ws.Range("E40").Value = "'2.76%"
$ws.Range("G40").Value = "'2"
$ws.Range("D41").Value = "'0.007314"
$ws.Range("E41").Value = "'-0.12%"
$ws.Range("G41").Value = "'2"
$ws.Range("D42").Value = "'0.008961"
$ws.Range("G42").Value = "'2"
$ws.Range("D43").Value = "'0.1329"
$ws.Range("E43").Value = "'2.23%"
$ws.Range("G43").Value = "'2"
$ws.Range("D44").Value = "'0.001947"
$ws.Range("E44").Value = "'-4.76%"
$ws.Range("G44").Value = "'2"
$ws.Range("D45").Value = "'0.009426"
$ws.Range("E45").Value = "'-1.41%"
$ws.Range("G45").Value = "'2"
$ws.Range("D46").Value = "'0.00006071"
$ws.Range("E46").Value = "'0.93%"
$ws.Range("G46").Value = "'2"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.12%"
$ws.Range("G47").Value = "'2"
$ws.Range("E48").Value = "'-0.03%"
$ws.Range("G48").Value = "'2"
$ws.Range("E49").Value = "'-3.43%"
$ws.Range("G49").Value = "'2"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.12%"
$ws.Range("G50").Value = "'2"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.12%"
$ws.Range("G51").Value = "'2"

Write-Host "Updated symbol list values successfully"
